# The template worksheet originally listed its sample-image rows starting
# at row 4 (Pol0_45_90_135, Img1_C1.tif, Img2_C1.tif in rows 4-6).
# This edit inserts two blank rows above that block, pushing the three
# existing rows down to rows 6-8, and updates the active selection to
# match (mirrors "Use assertThrown in ByExcelFinder" commit which added
# extra header/blank lines before the sample rows in the finder fixture).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new, completely blank rows before row 4 (shifting rows
# 4-6 down to 6-8). Clear() removes the default formatting that Insert()
# copies down from the row above, so the new rows stay empty/unstyled.
$ws.Rows.Item(4).Resize(2).Insert()
$ws.Rows.Item(4).Resize(2).Clear()

# Update the selection to match the new layout: active cell A6 with the
# whole moved block (A6:A8) selected.
[void]$ws.Range("A6:A8").Select()

Write-Output "Inserted 2 rows before row 4; sample rows now at 6-8."
